# Test data for Greece Market
# - Duplicate the "Croatia" worksheet (right after itself) to create a new
#   "Greece" worksheet, update its part-number cell, and make it the active
#   sheet/tab - mirroring how this data set is extended for each new market.

$wb = $excel.ActiveWorkbook

$croatia = $wb.Worksheets.Item("Croatia")

# Copy Croatia immediately after itself - Excel names the copy "Croatia (2)".
$croatia.Copy($null, $croatia)
$greece = $wb.Worksheets.Item("Croatia (2)")
$greece.Name = "Greece"

# Greece gets its own part number in B4 (new shared string).
$greece.Range("B4").Value = "NGC-4119/T3199"

# Restore Croatia's selection to "select all" and leave Greece selected at H18,
# matching the state Excel leaves behind after this kind of sheet duplication.
$croatia.Activate()
$croatia.Cells.Select()

$greece.Activate()
$greece.Range("H18").Select()
